$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2649.75
$ws.Range("J17").Value = 2699.6667
$ws.Range("L17").Value = 8099.000100000001
$ws.Range("N17").Value = -8435.000100000001

$ws.Range("H40").Value = 3698.4
$ws.Range("J40").Value = 3998.3333
$ws.Range("L40").Value = 3998.3333
$ws.Range("N40").Value = -4348.3333

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws.Range("H70").Value = 1899.1333
$ws.Range("I70").Value = 1048.25
$ws.Range("J70").Value = 2208.5454
$ws.Range("K70").Value = 3144.75
$ws.Range("L70").Value = 6625.6362
$ws.Range("M70").Value = -2874.75
$ws.Range("N70").Value = -7165.6362

$ws.Range("H73").Value = 1899.1333
$ws.Range("I73").Value = 1048.25
$ws.Range("J73").Value = 2208.5454
$ws.Range("K73").Value = 3144.75
$ws.Range("L73").Value = 6625.6362
$ws.Range("M73").Value = -2208.75
$ws.Range("N73").Value = -8497.636200000001

$ws.Range("H88").Value = 2615.2307
$ws.Range("I88").Value = 2722
$ws.Range("J88").Value = 2548.5
$ws.Range("K88").Value = 2722
$ws.Range("L88").Value = 2548.5
$ws.Range("M88").Value = -2316
$ws.Range("N88").Value = -3360.5

$ws.Range("H91").Value = 2615.2307
$ws.Range("I91").Value = 2722
$ws.Range("J91").Value = 2548.5
$ws.Range("K91").Value = 2722
$ws.Range("L91").Value = 2548.5
$ws.Range("M91").Value = -1318
$ws.Range("N91").Value = -5356.5

$ws.Range("H113").Value = 3690
$ws.Range("I113").Value = 3661.2856
$ws.Range("J113").Value = 3706.75
$ws.Range("K113").Value = 3661.2856
$ws.Range("L113").Value = 3706.75
$ws.Range("M113").Value = -407.2856000000002
$ws.Range("N113").Value = -10214.75

$ws.Range("H125").Value = 1069.5
$ws.Range("J125").Value = 1749.5
$ws.Range("L125").Value = 15745.5
$ws.Range("N125").Value = -20665.5

$ws.Range("H132").Value = 2487.4854
$ws.Range("I132").Value = 1966.2222
$ws.Range("K132").Value = 5898.6666
$ws.Range("M132").Value = -3368.6666

$ws.Range("H138").Value = 4410.9585
$ws.Range("I138").Value = 2465.4614
$ws.Range("J138").Value = 5510.587
$ws.Range("K138").Value = 7396.3842
$ws.Range("L138").Value = 16531.761
$ws.Range("M138").Value = -2256.3842
$ws.Range("N138").Value = -26811.761

$ws.Range("H141").Value = 8403.471
$ws.Range("I141").Value = 9123.933999999999
$ws.Range("K141").Value = 27371.802
$ws.Range("M141").Value = -22191.802

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 18314.285
$ws.Range("I97").Value = 35333.332
$ws.Range("K97").Value = 35333.332
$ws.Range("M97").Value = -34837.332

$ws.Range("H132").Value = 50841.285
$ws.Range("I132").Value = 64443.875
$ws.Range("K132").Value = 193331.625
$ws.Range("M132").Value = -190801.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 111899.22
$ws.Range("I22").Value = 500399.5
$ws.Range("J22").Value = 899.1429000000001
$ws.Range("K22").Value = 500399.5
$ws.Range("L22").Value = 899.1429000000001
$ws.Range("M22").Value = -500226.5
$ws.Range("N22").Value = -1245.1429

$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()

$ws.Range("H99").Value = 67530.94
$ws.Range("I99").Value = 116166.11
$ws.Range("K99").Value = 116166.11
$ws.Range("M99").Value = -114668.11

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2849.5
$ws.Range("I62").Value = 2800
$ws.Range("J62").Value = 2899
$ws.Range("K62").Value = 2800
$ws.Range("L62").Value = 2899
$ws.Range("M62").Value = -2176
$ws.Range("N62").Value = -4147

$ws.Range("H65").Value = 2849.5
$ws.Range("I65").Value = 2800
$ws.Range("J65").Value = 2899
$ws.Range("K65").Value = 14000
$ws.Range("L65").Value = 14495
$ws.Range("M65").Value = -10880
$ws.Range("N65").Value = -20735

$ws.Range("H86").Value = 17324.355
$ws.Range("I86").Value = 27522.117
$ws.Range("K86").Value = 27522.117
$ws.Range("M86").Value = -26399.117

$ws.Range("H89").Value = 17324.355
$ws.Range("I89").Value = 27522.117
$ws.Range("K89").Value = 137610.585
$ws.Range("M89").Value = -131994.585

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 6623.75
$ws.Range("I55").Value = 5000
$ws.Range("J55").Value = 6855.7144
$ws.Range("K55").Value = 15000
$ws.Range("L55").Value = 20567.1432
$ws.Range("M55").Value = -14823
$ws.Range("N55").Value = -20921.1432

$ws.Range("H80").Value = 33338000
$ws.Range("J80").Value = 10000
$ws.Range("L80").Value = 30000
$ws.Range("N80").Value = -31872

$ws.Range("H83").Value = 33338000
$ws.Range("J83").Value = 10000
$ws.Range("L83").Value = 90000
$ws.Range("N83").Value = -99360

$ws.Range("H108").Value = 3166.6667
$ws.Range("I108").Value = 3166.6667
$ws.Range("K108").Value = 9500.000100000001
$ws.Range("M108").Value = -6620.000100000001

$ws.Range("H115").Value = 994.5
$ws.Range("I115").Value = 994.5
$ws.Range("K115").Value = 2983.5
$ws.Range("M115").Value = -1808.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 784.7143
$ws.Range("I97").Value = 784.7143
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 784.7143
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -288.7143
$ws.Range("N97").ClearContents()

$ws.Range("H103").Value = 49667
$ws.Range("J103").Value = 49667
$ws.Range("L103").Value = 49667
$ws.Range("N103").Value = -52011

$ws.Range("H122").Value = 3362.9285
$ws.Range("I122").Value = 3362.9285
$ws.Range("K122").Value = 10088.7855
$ws.Range("M122").Value = -7638.7855

$ws.Range("H132").Value = 47708.87
$ws.Range("I132").Value = 55311.367
$ws.Range("J132").Value = 11597
$ws.Range("K132").Value = 165934.101
$ws.Range("L132").Value = 34791
$ws.Range("M132").Value = -163404.101
$ws.Range("N132").Value = -39851

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11683.866
$ws.Range("I7").Value = 17737.375
$ws.Range("J7").Value = 4765.5713
$ws.Range("K7").Value = 17737.375
$ws.Range("L7").Value = 4765.5713
$ws.Range("M7").Value = -17625.375
$ws.Range("N7").Value = -4989.5713

$ws.Range("H40").Value = 3210.375
$ws.Range("I40").Value = 2580.5
$ws.Range("K40").Value = 2580.5
$ws.Range("M40").Value = -2444.5

$ws.Range("H126").Value = 11683.866
$ws.Range("I126").Value = 17737.375
$ws.Range("J126").Value = 4765.5713
$ws.Range("K126").Value = 53212.125
$ws.Range("L126").Value = 14296.7139
$ws.Range("M126").Value = -50742.125
$ws.Range("N126").Value = -19236.7139

$ws.Range("H133").Value = 68000
$ws.Range("J133").Value = 68000
$ws.Range("L133").Value = 68000
$ws.Range("N133").Value = -73060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 20001736
$ws.Range("I96").Value = 2169.75
$ws.Range("K96").Value = 2169.75
$ws.Range("M96").Value = -796.75

$ws.Range("H100").Value = 793.6923
$ws.Range("I100").Value = 757.1818
$ws.Range("J100").Value = 994.5
$ws.Range("K100").Value = 1514.3636
$ws.Range("L100").Value = 1989
$ws.Range("M100").Value = -973.3635999999999
$ws.Range("N100").Value = -3071

$ws.Range("H107").Value = 937.625
$ws.Range("I107").Value = 1014.4286
$ws.Range("K107").Value = 3043.2858
$ws.Range("M107").Value = -1123.2858

$ws.Range("H122").Value = 524.9375
$ws.Range("I122").Value = 529.4167
$ws.Range("J122").Value = 511.5
$ws.Range("K122").Value = 1588.2501
$ws.Range("L122").Value = 1534.5
$ws.Range("M122").Value = 861.7499
$ws.Range("N122").Value = -6434.5

$ws.Range("H126").Value = 201037
$ws.Range("I126").Value = 212972.8
$ws.Range("K126").Value = 638918.3999999999
$ws.Range("M126").Value = -636448.3999999999

$ws.Range("H132").Value = 93290.66
$ws.Range("I132").Value = 109473.37
$ws.Range("K132").Value = 328420.11
$ws.Range("M132").Value = -325890.11

